# Apply "added new repair strategy CREATE_IF_MISSING" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Test-Payload (E20): add "unexpected": 17.5 field to the JSON sample
$ws.Range("E20").Value = "{`n     ""temperature"": 120.5,`n     ""unit"": ""Celsius"",`n     ""time"": ""2023-07-12T16:21:53.389+02:00"",`n     ""externalId"": ""berlin_01"",`n     ""unexpected"": 17.5`n}"

# Update Expected Result (G20): describe the new CREATE_IF_MISSING repair strategy
$ws.Range("G20").Value = "A measasurement should be created for the device berlin_01.`nThe fragment ""c8y_Fragment_to_remove"" is not included in the created measurement, as the repair strategy is ""REMOVE_IF_NULL"".`nIn addition the reapar strategy ""CREATE_IF_MISSING"" is used. Thjsi is required to map the node ""unexpected"" to the target fragment ""c8y_Unexpected"". This is created, due to the used reapir strategy."

# Excel auto row-height doesn't recompute in this headless runtime, so set explicitly
$ws.Rows.Item(20).RowHeight = 252

# Reflect the final cursor/scroll position left after editing
[void]$ws.Range("F20").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 2

$wb.Save()
